$d = $word.ActiveDocument

# 1. Update the date in the title line
$d.Content.Find.Execute("04.09.24", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "03.09.24", 2) | Out-Null

# 2. Replace the paper title paragraph (also removes the trailing <w:br/>)
$pTitle = $d.Paragraphs.Item(2)
$pTitle.Range.Text = "Smaller, Weaker, Yet Better: Training LLM Reasoners via Compute-Optimal Sampling"

# 3. Replace paragraph 3 body text
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "אחת הדרכים הדי מפתיעות לשיפור יכולות reasoning של מודלי שפה היא שיפור עצמי או self-improvement. בגדול עבור דאטהסט של שאלות ותשובות אנו מבקשים ממודל שפה לענות על התשובה ולספק הסבר. לאחר מכן מפלטרים את השרשראות reasoning שלא התכנסו לתשובה הרצויה. לאחר הפלטור מבצעים פיינטיון של המודל על הדאטהסט המפולטר. וכאמור באופן די מפתיע (לפחות אותי) הדבר אכן מוביל לשיפור יכולות reasoning של מודל שפה."

# 4. Replace paragraph 4 body text. The original run text ends with a trailing
#    space (xml:space="preserve"); use Find scoped to the paragraph so the
#    replacement text (no leading/trailing whitespace) does not inherit it.
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Find.Execute(
    "כמו שאתם בטח זוכרים LoRA היא משפחה (די גדולה שממשיכה לגדול) של שיטות מהמשפחה (גדולה עוד יותר) של שיטות חסכוניות פיינטיון של מודלי שפה ענקיים (או PEFT - Parameter Efficient Fine-Tuning). C ב-LoRA אנו מאמנים תוספת של משקלים לכל שכבה במקום לאמן את כל המשקלים במודל. כל תוספת כזו היא מטריצה בעלת רנק נמוך כלומר אפקטיבית מכילה מעט פרמטרים מאשר מטריצת המשקלים של השכבה. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "ואם יש בידינו מודל יותר חזק אז ניתן לבנות את הדאטהסט הזה באמצעותו ולעשות את הפיינטיון על הדאטה הנוצר באמצעותו בצורה דומה.",
    2) | Out-Null

# 5. Replace paragraph 5 body text
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = "אולם המאמר שואל שאלה די מעניינת: מה עדיף (מבחינת הביצועים), ליצור יחסית מעט דאטה עם מודל גדול וחזק או ליצור יחסית הרבה דאטה עם מודל קטן וחלש יותר. הרי יצירת דאטה עם מודל חזק היא יקרה יותר (מבחינת כמות ה-FLOPS הכוללת הנדרשת לכך) אבל מצד שני הדאטה שהוא יוצר הוא יותר איכותי."

# 6. Replace paragraph 6 body text. As with paragraph 4 above, the original
#    run carried xml:space="preserve"; scope the Find to the paragraph so it
#    is not inherited by the new (non-whitespace-padded) text.
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Find.Execute(
    "המאמר (וגם אני) חשבו על דרך אחרת של אופטימיזציה של LoRA. המחברים שואלים שאלה פשוטה - למה בנוסף לאימון של מטריצות התוספות לא נאמן את ה-importance שלה בכל שכבה. ה-importance במקרה הזה היא המקדם המכפיל את מטריצת התוספות לפני הוספתה מטריצת המשקלות המקורית במודל (שנותרת קבועה במהלך פיינטיון).  האלגוריתם המוצע עושה כמה איטרציות של משקלי ה-importance לעדכון אחד של משקלות התוספות. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "המחברים מציעים לבצע את ההשוואה של ״תפוזים לתפוזים״ - כלומר לקחת את הדאטה הנוצר עם מודל חזק ומודל חזק תחת אותו תקציב של FLOPS ולהשוות מה מהם מוביל לביצועים טובים יותר של המודל שעובר פיינטיון על הדאטה הזה.",
    2) | Out-Null

# 7. Replace paragraph 7 body text
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = "ויש תוצאות די מעניינות במאמר.."

# 8. Replace paragraph 8 (link paragraph): drop the <w:br/> and shrink to "."
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = "."

# 9. Append a new paragraph with the updated link
$d.Paragraphs.Item(8).Range.InsertParagraphAfter()
$d.Paragraphs.Item(9).Range.Text = "https://arxiv.org/pdf/2408.16737"
